$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "31/12/2006"
$ws.Range("C2").Value = 19.0601393249186

$ws.Range("B3").Value = "31/12/2010"
$ws.Range("C3").Value = 14.9909728655381

$ws.Range("B4").Value = "31/12/2014"
$ws.Range("C4").Value = 13.4671494565605

$ws.Range("B5").Value = "31/12/2018"
$ws.Range("C5").Value = 12.4522364098828

$ws.Range("B6").Value = "31/12/2022"
$ws.Range("C6").Value = 12.4213105471151

$ws.Range("B7").Value = "31/12/2024"
$ws.Range("C7").Value = 12.2488303378121

$ws.Range("B8").Value = "31/12/2006"
$ws.Range("C8").Value = 24.5300653231751

$ws.Range("B9").Value = "31/12/2010"
$ws.Range("C9").Value = 17.6707457583977

$ws.Range("B10").Value = "31/12/2014"
$ws.Range("C10").Value = 15.5352034014426

$ws.Range("B11").Value = "31/12/2018"
$ws.Range("C11").Value = 13.9909085963737

$ws.Range("B12").Value = "31/12/2022"
$ws.Range("C12").Value = 13.6786420123094

$ws.Range("B13").Value = "31/12/2024"
$ws.Range("C13").Value = 13.4480186995168

$ws.Range("B14").Value = "31/12/2006"
$ws.Range("C14").Value = 26.0096105678633

$ws.Range("B15").Value = "31/12/2010"
$ws.Range("C15").Value = 16.332031721116

$ws.Range("B16").Value = "31/12/2014"
$ws.Range("C16").Value = 16.2918843775623

$ws.Range("B17").Value = "31/12/2018"
$ws.Range("C17").Value = 17.6566831234739

$ws.Range("B18").Value = "31/12/2022"
$ws.Range("C18").Value = 16.5584731824403

$ws.Range("B19").Value = "31/12/2024"
$ws.Range("C19").Value = 18.2449622065689
